$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Update column C (linea) values
$ws.Range("C2").Value = 0.43
$ws.Range("C3").Value = 0.7
$ws.Range("C4").Value = 0.55000000000000004
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C9").Value = 2
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 0

# Update column F (vestizione) values that also changed
$ws.Range("F5").Value = 10.18
$ws.Range("F6").Value = 6.33

# Update the active selection to match the saved view state
$ws.Range("L13").Select()
